$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2021 column (R) is being removed from the sheet entirely - delete the
# whole column and shift everything left (also removes the R4:R14 values).
$ws.Range("R1:R14").EntireColumn.Delete()

# Update the selection to match the new state of the sheet.
$ws.Range("N19").Select()
